$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the old test-step identifiers in column A (prefix with "old_"),
# preserving the order in which the "old_" strings were introduced so the
# shared-string table is rebuilt in the same order as the target workbook.
$ws.Range("A12").Value = "old_Choose Silver"
$ws.Range("A13").Value = "old_Choose Gold"
$ws.Range("A14").Value = "old_Choose Platinum"
$ws.Range("A15").Value = "old_Choose Ultimate"
$ws.Range("A7").Value  = "old_102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageVehicleData"
$ws.Range("A8").Value  = "old_102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageInsurantData"
$ws.Range("A9").Value  = "old_102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageProductData"
$ws.Range("A10").Value = "old_102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageSendQuote"

# Highlight the renamed ("now inactive") rows with an orange fill so they are
# visually distinguishable from the active steps.
$ws.Range("A7").Interior.Color = 49407
$ws.Range("A8").Interior.Color = 49407
$ws.Range("A9").Interior.Color = 49407
$ws.Range("A10").Interior.Color = 49407
$ws.Range("A12").Interior.Color = 49407
$ws.Range("A13").Interior.Color = 49407
$ws.Range("A14").Interior.Color = 49407
$ws.Range("A15").Interior.Color = 49407

# Move the active selection to A7, matching the saved cursor position.
$ws.Range("A7").Select()
